$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two placeholder rows ("Personne" and the all-zero "Kazsc" row)
# that sat at rows 2 and 3. This shifts every following player's row up by
# two positions (e.g. former row 4 "ylarabka" becomes row 2, etc.) and the
# sheet's used range shrinks from A1:S13 down to A1:S11.
$ws.Rows("2:3").Delete()

# After the shift, rows 4 (tomlora), 5 (chatobogan) and 10 (namiyeon) also
# received updated statistics, so refresh those rows explicitly, cell by cell.
$row4 = @("tomlora", 1, 1, 123, 105, 61.22950000000002, 2507, 1071, 346, 338, 25067, 1025, 748, 927, 8.333333333333334, 6.08130081300813, 7.536585365853658, 20.38, 29.87)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

$row5 = @("chatobogan", 0, 0, 46, 38, 24.10116666666666, 2480, 1201, 341, 329, 4476, 244, 185, 539, 5.304347826086956, 4.021739130434782, 11.71739130434783, 53.91, 31.44)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

$row10 = @("namiyeon", 0, 0, 62, 6, 31.29, 4969, 2837, 613, 672, 1203, 150, 186, 1028, 2.419354838709677, 3, 16.58064516129032, 80.15000000000001, 30.28)
for ($i = 0; $i -lt $row10.Length; $i++) {
    $ws.Cells.Item(10, $i + 1).Value = $row10[$i]
}
